$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H5").Value = 250000080
$ws.Range("I5").Value = 500000100
$ws.Range("J5").Value = 83.5
$ws.Range("K5").Value = 500000100
$ws.Range("L5").Value = 83.5
$ws.Range("M5").Value = -499999985
$ws.Range("N5").Value = -313.5
$ws.Range("H15").Value = 1181.919
$ws.Range("I15").Value = 1181.919
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 3545.757000000001
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -3376.757000000001
$ws.Range("H48").Value = 1700
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 1700
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 5100
$ws.Range("N48").Value = -5684
$ws.Range("H56").Value = 1700
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 1700
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 5100
$ws.Range("N56").Value = -6168
$ws.Range("H74").Value = 4499.3335
$ws.Range("I74").Value = 4499
$ws.Range("H77").Value = 4499.3335
$ws.Range("I77").Value = 4499
$ws.Range("H101").Value = 595.1111
$ws.Range("I101").Value = 746.1667
$ws.Range("J101").Value = 293
$ws.Range("K101").Value = 2238.5001
$ws.Range("L101").Value = 879
$ws.Range("M101").Value = -616.5001000000002
$ws.Range("N101").Value = -4123
$ws.Range("H106").Value = 38399.465
$ws.Range("I106").Value = 39153.23
$ws.Range("J106").Value = 33500
$ws.Range("K106").Value = 39153.23
$ws.Range("L106").Value = 33500
$ws.Range("M106").Value = -38522.23
$ws.Range("N106").Value = -34762
$ws.Range("H137").Value = 1804.0714
$ws.Range("I137").Value = 1717.5555
$ws.Range("J137").Value = 1959.8
$ws.Range("K137").Value = 5152.666499999999
$ws.Range("L137").Value = 5879.4
$ws.Range("M137").Value = -2602.666499999999
$ws.Range("N137").Value = -10979.4

# ---- Sheet: ARM ----
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H30").Value = 1941.4286
$ws.Range("I30").Value = 1795
$ws.Range("J30").Value = 2000
$ws.Range("K30").Value = 1795
$ws.Range("L30").Value = 2000
$ws.Range("M30").Value = -1645
$ws.Range("N30").Value = -2300
$ws.Range("H32").Value = 3800.0925
$ws.Range("I32").Value = 2404.18
$ws.Range("J32").Value = 21249
$ws.Range("K32").Value = 2404.18
$ws.Range("L32").Value = 21249
$ws.Range("M32").Value = -2117.18
$ws.Range("N32").Value = -21823
$ws.Range("H61").Value = 1416.909
$ws.Range("I61").Value = 1318.6
$ws.Range("J61").Value = 2400
$ws.Range("K61").Value = 1318.6
$ws.Range("L61").Value = 2400
$ws.Range("M61").Value = -1106.6
$ws.Range("N61").Value = -2824
$ws.Range("H132").Value = 1294.8422
$ws.Range("I132").Value = 1247.6666
$ws.Range("J132").Value = 1471.75
$ws.Range("K132").Value = 3742.9998
$ws.Range("L132").Value = 4415.25
$ws.Range("M132").Value = -1212.9998
$ws.Range("N132").Value = -9475.25
$ws.Range("H136").Value = 1416.909
$ws.Range("I136").Value = 1318.6
$ws.Range("J136").Value = 2400
$ws.Range("K136").Value = 3955.8
$ws.Range("L136").Value = 7200
$ws.Range("M136").Value = -1405.8
$ws.Range("N136").Value = -12300

# ---- Sheet: BSM ----
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H86").Value = 2036.5
$ws.Range("I86").Value = 2665.3333
$ws.Range("J86").Value = 150
$ws.Range("K86").Value = 2665.3333
$ws.Range("L86").Value = 150
$ws.Range("M86").Value = -1542.3333
$ws.Range("N86").Value = -2396
$ws.Range("H89").Value = 2036.5
$ws.Range("I89").Value = 2665.3333
$ws.Range("J89").Value = 150
$ws.Range("K89").Value = 13326.6665
$ws.Range("L89").Value = 750
$ws.Range("M89").Value = -7710.666499999999
$ws.Range("N89").Value = -11982
$ws.Range("H129").Value = 70780
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 70780
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 70780
$ws.Range("N129").Value = -80780

# ---- Sheet: CRP ----
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H7").Value = 105.73684
$ws.Range("I7").Value = 50
$ws.Range("J7").Value = 261.8
$ws.Range("K7").Value = 50
$ws.Range("L7").Value = 261.8
$ws.Range("M7").Value = 63
$ws.Range("N7").Value = -487.8
$ws.Range("H31").Value = 3392
$ws.Range("I31").Value = 2994.4167
$ws.Range("J31").Value = 5777.5
$ws.Range("K31").Value = 2994.4167
$ws.Range("L31").Value = 5777.5
$ws.Range("M31").Value = -2699.4167
$ws.Range("N31").Value = -6367.5
$ws.Range("H34").Value = 3392
$ws.Range("I34").Value = 2994.4167
$ws.Range("J34").Value = 5777.5
$ws.Range("K34").Value = 2994.4167
$ws.Range("L34").Value = 5777.5
$ws.Range("M34").Value = -2792.4167
$ws.Range("N34").Value = -6181.5
$ws.Range("H58").Value = 1910.3235
$ws.Range("I58").Value = 1120.1724
$ws.Range("J58").Value = 6493.2
$ws.Range("K58").Value = 1120.1724
$ws.Range("L58").Value = 6493.2
$ws.Range("M58").Value = -917.1723999999999
$ws.Range("N58").Value = -6899.2
$ws.Range("H99").Value = 12211.667
$ws.Range("I99").Value = 9107.223
$ws.Range("J99").Value = 14540
$ws.Range("K99").Value = 9107.223
$ws.Range("L99").Value = 14540
$ws.Range("M99").Value = -7609.223
$ws.Range("N99").Value = -17536
$ws.Range("H105").Value = 1226.1111
$ws.Range("I105").Value = 1244.375
$ws.Range("J105").Value = 1080
$ws.Range("K105").Value = 1244.375
$ws.Range("L105").Value = 1080
$ws.Range("M105").Value = 502.625
$ws.Range("N105").Value = -4574
$ws.Range("H107").Value = 849.4167
$ws.Range("I107").Value = 495.75
$ws.Range("J107").Value = 1556.75
$ws.Range("K107").Value = 495.75
$ws.Range("L107").Value = 1556.75
$ws.Range("M107").Value = 1424.25
$ws.Range("N107").Value = -5396.75
$ws.Range("H126").Value = 12211.667
$ws.Range("I126").Value = 9107.223
$ws.Range("J126").Value = 14540
$ws.Range("K126").Value = 27321.669
$ws.Range("L126").Value = 43620
$ws.Range("M126").Value = -24851.669
$ws.Range("N126").Value = -48560
$ws.Range("H136").Value = 1910.3235
$ws.Range("I136").Value = 1120.1724
$ws.Range("J136").Value = 6493.2
$ws.Range("K136").Value = 3360.5172
$ws.Range("L136").Value = 19479.6
$ws.Range("M136").Value = -810.5171999999998
$ws.Range("N136").Value = -24579.6

# ---- Sheet: CUL ----
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H4").Value = 3797189
$ws.Range("I4").Value = 4673304.5
$ws.Range("J4").Value = 687.8333
$ws.Range("K4").Value = 14019913.5
$ws.Range("L4").Value = 2063.4999
$ws.Range("M4").Value = -14019801.5
$ws.Range("N4").Value = -2287.4999
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H121").Value = 856
$ws.Range("I121").Value = 800
$ws.Range("J121").Value = 863
$ws.Range("K121").Value = 2400
$ws.Range("L121").Value = 2589
$ws.Range("M121").Value = -1090
$ws.Range("N121").Value = -5209

# ---- Sheet: GSM ----
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H2").Value = 753.9
$ws.Range("I2").Value = 80
$ws.Range("J2").Value = 922.375
$ws.Range("K2").Value = 80
$ws.Range("L2").Value = 922.375
$ws.Range("M2").Value = 33
$ws.Range("N2").Value = -1148.375
$ws.Range("H70").Value = 5831.25
$ws.Range("I70").Value = 5162.875
$ws.Range("J70").Value = 6499.625
$ws.Range("K70").Value = 5162.875
$ws.Range("L70").Value = 6499.625
$ws.Range("M70").Value = -4892.875
$ws.Range("N70").Value = -7039.625
$ws.Range("H73").Value = 5831.25
$ws.Range("I73").Value = 5162.875
$ws.Range("J73").Value = 6499.625
$ws.Range("K73").Value = 5162.875
$ws.Range("L73").Value = 6499.625
$ws.Range("M73").Value = -4226.875
$ws.Range("N73").Value = -8371.625
$ws.Range("H101").Value = 100100
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 100100
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 100100
$ws.Range("N101").Value = -106590
$ws.Range("H102").Value = 2216
$ws.Range("I102").Value = 1961.1428
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 1961.1428
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -339.1428000000001
$ws.Range("N102").Value = -7244
$ws.Range("H113").Value = 3110.75
$ws.Range("I113").Value = 2221.5
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2221.5
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -51.5
$ws.Range("N113").Value = -8340
$ws.Range("H132").Value = 2725.8
$ws.Range("I132").Value = 2562
$ws.Range("J132").Value = 4200
$ws.Range("K132").Value = 7686
$ws.Range("L132").Value = 12600
$ws.Range("M132").Value = -5156
$ws.Range("N132").Value = -17660

# ---- Sheet: LTW ----
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H5").Value = 22270.334
$ws.Range("I5").Value = 1800
$ws.Range("J5").Value = 32505.5
$ws.Range("K5").Value = 1800
$ws.Range("L5").Value = 32505.5
$ws.Range("M5").Value = -1687
$ws.Range("N5").Value = -32731.5
$ws.Range("H22").Value = 3189.4443
$ws.Range("I22").Value = 1797.2142
$ws.Range("J22").Value = 8062.25
$ws.Range("K22").Value = 1797.2142
$ws.Range("L22").Value = 8062.25
$ws.Range("M22").Value = -1502.2142
$ws.Range("N22").Value = -8652.25
$ws.Range("H27").Value = 3189.4443
$ws.Range("I27").Value = 1797.2142
$ws.Range("J27").Value = 8062.25
$ws.Range("K27").Value = 1797.2142
$ws.Range("L27").Value = 8062.25
$ws.Range("M27").Value = -1690.2142
$ws.Range("N27").Value = -8276.25
$ws.Range("H68").Value = 5501
$ws.Range("I68").Value = 5333.3335
$ws.Range("J68").Value = 5668.6665
$ws.Range("K68").Value = 5333.3335
$ws.Range("L68").Value = 5668.6665
$ws.Range("M68").Value = -4584.3335
$ws.Range("N68").Value = -7166.6665
$ws.Range("H71").Value = 5501
$ws.Range("I71").Value = 5333.3335
$ws.Range("J71").Value = 5668.6665
$ws.Range("K71").Value = 26666.6675
$ws.Range("L71").Value = 28343.3325
$ws.Range("M71").Value = -22922.6675
$ws.Range("N71").Value = -35831.3325
$ws.Range("H122").Value = 14249.5
$ws.Range("I122").Value = 14249.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 42748.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -40298.5
$ws.Range("H132").Value = 2637.2856
$ws.Range("I132").Value = 2475.8235
$ws.Range("J132").Value = 3323.5
$ws.Range("K132").Value = 7427.470499999999
$ws.Range("L132").Value = 9970.5
$ws.Range("M132").Value = -4897.470499999999
$ws.Range("N132").Value = -15030.5
$ws.Range("H136").Value = 2764.111
$ws.Range("I136").Value = 2764.111
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8292.332999999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5742.332999999999

# ---- Sheet: WVR ----
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H21").Value = 12979.25
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 12979.25
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 12979.25
$ws.Range("N21").Value = -13449.25
$ws.Range("H25").Value = 10000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 10000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 10000
$ws.Range("N25").Value = -10586
$ws.Range("H35").Value = 12979.25
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 12979.25
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 12979.25
$ws.Range("N35").Value = -13559.25
$ws.Range("H37").Value = 40029
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 40029
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 40029
$ws.Range("N37").Value = -40435
$ws.Range("H123").Value = 16666.334
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 16666.334
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 16666.334
$ws.Range("N123").Value = -26466.334
$ws.Range("H132").Value = 50469.42
$ws.Range("I132").Value = 73090.08
$ws.Range("J132").Value = 1458
$ws.Range("K132").Value = 219270.24
$ws.Range("L132").Value = 219270.24
$ws.Range("M132").Value = -216740.24
$ws.Range("N132").Value = -9434
